# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the crafting-class
# sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 103.40909
$ws.Range("I11").Value = 103.40909
$ws.Range("K11").Value = 103.40909
$ws.Range("M11").Value = 36.59090999999999
$ws.Range("H86").Value = 2923.4666
$ws.Range("I86").Value = 2661
$ws.Range("K86").Value = 2661
$ws.Range("M86").Value = -1538
$ws.Range("H89").Value = 2923.4666
$ws.Range("I89").Value = 2661
$ws.Range("K89").Value = 13305
$ws.Range("M89").Value = -7689

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2726.8462
$ws.Range("I2").Value = 1135.5714
$ws.Range("K2").Value = 1135.5714
$ws.Range("M2").Value = -1022.5714
$ws.Range("H4").Value = 264.6
$ws.Range("J4").Value = 246.25
$ws.Range("L4").Value = 246.25
$ws.Range("N4").Value = -478.25
$ws.Range("H5").Value = 133.16667
$ws.Range("J5").Value = 99
$ws.Range("L5").Value = 99
$ws.Range("N5").Value = -323
$ws.Range("H45").Value = 8490.286
$ws.Range("I45").Value = 20995.8
$ws.Range("J45").Value = 1542.7778
$ws.Range("K45").Value = 20995.8
$ws.Range("L45").Value = 1542.7778
$ws.Range("M45").Value = -20618.8
$ws.Range("N45").Value = -2296.7778
$ws.Range("H61").Value = 5323.5884
$ws.Range("I61").Value = 5677.467
$ws.Range("J61").Value = 2669.5
$ws.Range("K61").Value = 5677.467
$ws.Range("L61").Value = 2669.5
$ws.Range("M61").Value = -5465.467
$ws.Range("N61").Value = -3093.5
$ws.Range("H88").Value = 3727.5454
$ws.Range("I88").Value = 3000.2
$ws.Range("K88").Value = 3000.2
$ws.Range("M88").Value = -2594.2
$ws.Range("H91").Value = 3727.5454
$ws.Range("I91").Value = 3000.2
$ws.Range("K91").Value = 3000.2
$ws.Range("M91").Value = -1596.2
$ws.Range("H116").Value = 2726.8462
$ws.Range("I116").Value = 1135.5714
$ws.Range("K116").Value = 1135.5714
$ws.Range("M116").Value = 1158.4286
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 8700
$ws.Range("N122").Value = -13600
$ws.Range("H132").Value = 2007.56
$ws.Range("I132").Value = 1964.7391
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5894.2173
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3364.2173
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 5323.5884
$ws.Range("I136").Value = 5677.467
$ws.Range("J136").Value = 2669.5
$ws.Range("K136").Value = 17032.401
$ws.Range("L136").Value = 8008.5
$ws.Range("M136").Value = -14482.401
$ws.Range("N136").Value = -13108.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2726.8462
$ws.Range("I3").Value = 1135.5714
$ws.Range("K3").Value = 1135.5714
$ws.Range("M3").Value = -1021.5714
$ws.Range("H4").Value = 133.16667
$ws.Range("J4").Value = 99
$ws.Range("L4").Value = 99
$ws.Range("N4").Value = -329
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H105").Value = 4151.4653
$ws.Range("J105").Value = 1811
$ws.Range("L105").Value = 1811
$ws.Range("N105").Value = -5305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 625.8125
$ws.Range("I5").Value = 401.44446
$ws.Range("J5").Value = 914.2857
$ws.Range("K5").Value = 401.44446
$ws.Range("L5").Value = 914.2857
$ws.Range("M5").Value = -289.44446
$ws.Range("N5").Value = -1138.2857
$ws.Range("H25").Value = 2284.9
$ws.Range("I25").Value = 205.55556
$ws.Range("K25").Value = 205.55556
$ws.Range("M25").Value = -31.55556000000001
$ws.Range("H31").Value = 6583.4116
$ws.Range("J31").Value = 10999.125
$ws.Range("L31").Value = 10999.125
$ws.Range("N31").Value = -11589.125
$ws.Range("H34").Value = 6583.4116
$ws.Range("J34").Value = 10999.125
$ws.Range("L34").Value = 10999.125
$ws.Range("N34").Value = -11403.125
$ws.Range("H37").Value = 6642.857
$ws.Range("H127").Value = 49916.668
$ws.Range("I127").Value = 30000
$ws.Range("J127").Value = 59875
$ws.Range("K127").Value = 30000
$ws.Range("L127").Value = 59875
$ws.Range("M127").Value = -25040
$ws.Range("N127").Value = -69795
$ws.Range("I132").Value = 1649.4546
$ws.Range("J132").Value = 2065
$ws.Range("K132").Value = 4948.3638
$ws.Range("L132").Value = 6195
$ws.Range("M132").Value = -2418.3638
$ws.Range("N132").Value = -11255
$ws.Range("H137").Value = 71417.836
$ws.Range("J137").Value = 71417.836
$ws.Range("L137").Value = 71417.836
$ws.Range("N137").Value = -81617.836

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 395.1
$ws.Range("I11").Value = 300.16666
$ws.Range("K11").Value = 900.4999799999999
$ws.Range("M11").Value = -760.4999799999999
$ws.Range("H17").Value = 9182022
$ws.Range("I17").Value = 14285854
$ws.Range("J17").Value = 250314.75
$ws.Range("K17").Value = 42857562
$ws.Range("L17").Value = 750944.25
$ws.Range("M17").Value = -42857393
$ws.Range("N17").Value = -751282.25
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H46").Value = 17940298
$ws.Range("J46").Value = 5001131
$ws.Range("L46").Value = 15003393
$ws.Range("N46").Value = -15003575
$ws.Range("H55").Value = 625493
$ws.Range("I55").Value = 625493
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1876479
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1876302
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H132").Value = 2366.7144
$ws.Range("I132").Value = 1627
$ws.Range("K132").Value = 14643
$ws.Range("M132").Value = -12113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5603.222
$ws.Range("I2").Value = 9141.362999999999
$ws.Range("J2").Value = 43.285713
$ws.Range("K2").Value = 9141.362999999999
$ws.Range("L2").Value = 43.285713
$ws.Range("M2").Value = -9028.362999999999
$ws.Range("N2").Value = -269.285713
$ws.Range("H49").Value = 26249.75
$ws.Range("J49").Value = 26249.75
$ws.Range("L49").Value = 26249.75
$ws.Range("N49").Value = -26617.75
$ws.Range("H122").Value = 16593.234
$ws.Range("I122").Value = 17472.4
$ws.Range("K122").Value = 52417.2
$ws.Range("M122").Value = -49967.2
$ws.Range("H123").Value = 37498.75
$ws.Range("J123").Value = 37498.75
$ws.Range("L123").Value = 37498.75
$ws.Range("N123").Value = -42398.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 6248.5
$ws.Range("J11").Value = 6664.6665
$ws.Range("L11").Value = 6664.6665
$ws.Range("N11").Value = -6944.6665
$ws.Range("H17").Value = 1611.8889
$ws.Range("I17").Value = 1102.6086
$ws.Range("J17").Value = 2512.923
$ws.Range("K17").Value = 1102.6086
$ws.Range("L17").Value = 2512.923
$ws.Range("M17").Value = -932.6086
$ws.Range("N17").Value = -2852.923
$ws.Range("H18").Value = 6970
$ws.Range("I18").Value = 6940
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 6940
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = -6768
$ws.Range("N18").Value = -7344
$ws.Range("H25").Value = 9189
$ws.Range("I25").Value = 7000
$ws.Range("J25").Value = 9736.25
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 9736.25
$ws.Range("M25").Value = -6770
$ws.Range("N25").Value = -10196.25
$ws.Range("H46").Value = 2891.8462
$ws.Range("I46").Value = 1608.3334
$ws.Range("K46").Value = 1608.3334
$ws.Range("M46").Value = -1420.3334
$ws.Range("H62").Value = 49900
$ws.Range("J62").Value = 49900
$ws.Range("L62").Value = 49900
$ws.Range("N62").Value = -51148
$ws.Range("H65").Value = 49900
$ws.Range("J65").Value = 49900
$ws.Range("L65").Value = 149700
$ws.Range("N65").Value = -155940
$ws.Range("H88").Value = 31171
$ws.Range("I88").Value = 31171
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 31171
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("M88").Value = -30743
$ws.Range("H91").Value = 31171
$ws.Range("I91").Value = 31171
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 31171
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("M91").Value = -29689

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5175400.5
$ws.Range("I2").Value = 5175400.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5175400.5
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -5175288.5
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0
